# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets of the 广州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 192
$ws1.Range("F12").Value = 64
$ws1.Range("F15").Value = 1945
$ws1.Range("F16").Value = 447
$ws1.Range("F17").Value = 6540

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 510
$ws2.Range("F11").Value = 22
$ws2.Range("F12").Value = 112

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5445
$ws3.Range("F3").Value = 375
$ws3.Range("F4").Value = 369

# Sheet "全部类型" (All Types) - aggregated view of the above sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5445
$ws4.Range("F4").Value = 375
$ws4.Range("F5").Value = 369
$ws4.Range("F10").Value = 510
$ws4.Range("F22").Value = 192
$ws4.Range("F24").Value = 64
$ws4.Range("F25").Value = 22
$ws4.Range("F28").Value = 112
$ws4.Range("F29").Value = 1945
$ws4.Range("F30").Value = 447
$ws4.Range("F31").Value = 6540

$wb.Save()
